$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference date (column G) for every data row (2-274) from
# 45512 (2024-08-08) to 45513 (2024-08-09) - this report was re-run a day later.
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45513
}

# A handful of accounts now carry a projected value (column D) and the
# total (column H) reflects Vl. Projetado + Saldo Previsto for those rows.
$ws.Cells.Item(5, 4).Value = -10505.04
$ws.Cells.Item(5, 8).Value = -2994.88

$ws.Cells.Item(43, 4).Value = -58.34
$ws.Cells.Item(43, 8).Value = 647.54999999999995

$ws.Cells.Item(60, 4).Value = -16962.419999999998
$ws.Cells.Item(60, 8).Value = -140.76

$ws.Cells.Item(271, 4).Value = -9198.48
$ws.Cells.Item(271, 8).Value = -1242.02

# Sheet name reflects the new export timestamp.
$ws.Name = "IClientBalance-20240809-090208-"
